# Add season-record columns (Wins/Losses/Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new labels in AD1:AF1
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the existing header formatting (bold font, thin border, centered)
# by copying the format from an existing header cell.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Season record values for every player row (2 through 48).
$lastRow = 48
$ws.Range("AD2:AD" + $lastRow).Value = 82
$ws.Range("AE2:AE" + $lastRow).Value = 80
$ws.Range("AF2:AF" + $lastRow).Value = 0
